# Add data for 2022-04-02 refresh: extend the "through March 24" window to
# "through March 25" and bump the corresponding carjacking counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet tab name + the "March 2022" column header text both mention the
# as-of date; update both to the new cutoff.
$ws.Name = "Through 2022-03-25"
$ws.Range("B1").Value = "March 2022 (through March 25)"

# Per-neighborhood cell updates (row = neighborhood, column = month bucket).
$ws.Range("H3").Value = 4     # Austin
$ws.Range("Q4").Value = 3     # North Lawndale
$ws.Range("H5").Value = 5     # Garfield Park
$ws.Range("K5").Value = 1     # Garfield Park (new value)
$ws.Range("N11").Value = 6    # Englewood
$ws.Range("E17").Value = 2    # Auburn Gresham
$ws.Range("N17").Value = 1    # Auburn Gresham (new value)
$ws.Range("K19").Value = 2    # Lincoln Park
$ws.Range("N21").Value = 1    # West Pullman (new value)
$ws.Range("E22").Value = 1    # Bridgeport (new value)
$ws.Range("B26").Value = 2    # Chatham
$ws.Range("E26").Value = 4    # Chatham
$ws.Range("N26").Value = 2    # Chatham
$ws.Range("B27").Value = 2    # Calumet Heights
$ws.Range("T32").Value = 3    # New City
$ws.Range("E70").Value = 1    # Lincoln Square (new value)
$ws.Range("B88").Value = 2    # Uptown
